# Revising tables to improve readability
# Replace the middle-dot decimal separator ("·", U+00B7) with a standard
# period ("." ) in the numeric values of the data rows of Table 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the table shape on the slide (rather than assuming a fixed index).
$tblShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tblShape = $candidate
        break
    }
}
$tbl = $tblShape.Table

# Map of (row, column) -> corrected text, using 1-based COM indices.
# Row 3 = "1x GDP/capita", Row 4 = "3x GDP/capita" data rows.
$updates = @(
    @{ Row = 3; Col = 2;  Text = "1.42" },
    @{ Row = 3; Col = 3;  Text = "1.06" },
    @{ Row = 3; Col = 4;  Text = "1.85" },
    @{ Row = 3; Col = 6;  Text = "3.64" },
    @{ Row = 3; Col = 7;  Text = "2.25" },
    @{ Row = 3; Col = 8;  Text = "5.73" },
    @{ Row = 3; Col = 10; Text = "4.74" },
    @{ Row = 3; Col = 11; Text = "3.14" },
    @{ Row = 3; Col = 12; Text = "6.90 " },
    @{ Row = 4; Col = 2;  Text = "4.27" },
    @{ Row = 4; Col = 3;  Text = "3.17" },
    @{ Row = 4; Col = 4;  Text = "5.55" },
    @{ Row = 4; Col = 6;  Text = "10.93" },
    @{ Row = 4; Col = 7;  Text = "6.76" },
    @{ Row = 4; Col = 8;  Text = "17.20 " },
    @{ Row = 4; Col = 10; Text = "14.22" },
    @{ Row = 4; Col = 11; Text = "9.42" },
    @{ Row = 4; Col = 12; Text = "20.71" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cell.Shape.TextFrame.TextRange.Text = $u.Text
}
